# Weekly fruit/vegetable price update: a new daily record is inserted as
# row 12 (pushing all subsequent records down by one row), matching the
# logic used by the daily "subconjuntos" consolidation sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; everything below (old rows 12..85)
# shifts down to 13..86, and the sheet's used range grows to R86.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new record's data.
$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(12, 3).Value = "Metropolitana"
$ws.Cells.Item(12, 4).Value = 44532
$ws.Cells.Item(12, 5).Value = 13
$ws.Cells.Item(12, 6).Value = 100112022
$ws.Cells.Item(12, 7).Value = "Arveja Verde"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 43
$ws.Cells.Item(12, 11).Value = 16000
$ws.Cells.Item(12, 12).Value = 17000
$ws.Cells.Item(12, 13).Value = 16512
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Región del Maule"
$ws.Cells.Item(12, 16).Value = 660
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
